$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: sofiya -> sofiyaansari, 342524 -> 12345
$ws.Range("A2").Value = "sofiyaansari"
$ws.Range("B2").Value = 12345

# Update row 3: salamna/09654/m -> Saniya khan/12346/f
$ws.Range("A3").Value = "Saniya khan"
$ws.Range("B3").Value = 12346
$ws.Range("D3").Value = "f"

# Update row 4: salamna/09654/m -> sanakhan/12347/f
$ws.Range("A4").Value = "sanakhan"
$ws.Range("B4").Value = 12347
$ws.Range("D4").Value = "f"

# Update row 5: salamna/09654 -> salamna/96540 (name unchanged)
$ws.Range("B5").Value = 96540

# Update row 6: salamna/09654 -> salamnakhan/96542
$ws.Range("A6").Value = "salamnakhan"
$ws.Range("B6").Value = 96542

# Update row 7: sofiya/342524 -> sofiyaalam/34200
$ws.Range("A7").Value = "sofiyaalam"
$ws.Range("B7").Value = 34200

# Update row 8: salamna/09654 -> Maaz/9054
$ws.Range("A8").Value = "Maaz"
$ws.Range("B8").Value = 9054

# Update row 9: sofiya/342524 -> Vidhi/302504
$ws.Range("A9").Value = "Vidhi"
$ws.Range("B9").Value = 302504

# Update row 10: 98765 -> 987650 (name/designation/gender unchanged)
$ws.Range("B10").Value = 987650

# Update row 11: 543262 -> 5432020 (name/designation/gender unchanged)
$ws.Range("B11").Value = 5432020

# Update row 12: 345672 -> 305672 (name/designation/gender unchanged)
$ws.Range("B12").Value = 305672

$ws.Range("C15").Select()
